# Replace the multiplication problem/answer strings throughout the document.
# Each old value is unique in the document, so a simple Find/Replace (ReplaceAll)
# on the whole document content is safe and will not create unwanted side effects.

$d = $word.ActiveDocument

$replacements = @(
    @("858×7=6006", "111×6=666"),
    @("472×8=3776", "534×7=3738"),
    @("449×6=2694", "468×3=1404"),
    @("116×3=348",  "390×2=780"),
    @("623×5=3115", "891×8=7128"),
    @("681×5=3405", "292×9=2628"),
    @("529×2=1058", "973×5=4865"),
    @("580×4=2320", "379×2=758"),
    @("475×6=2850", "380×7=2660"),
    @("777×8=6216", "704×7=4928"),
    @("871×2=1742", "470×4=1880"),
    @("825×5=4125", "545×5=2725"),
    @("871×9=7839", "755×6=4530"),
    @("217×6=1302", "389×9=3501"),
    @("329×6=1974", "801×8=6408"),
    @("317×4=1268", "337×7=2359"),
    @("734×9=6606", "345×6=2070"),
    @("982×8=7856", "486×8=3888"),
    @("828×2=1656", "891×9=8019"),
    @("301×7=2107", "446×7=3122"),
    @("297×4=1188", "696×4=2784"),
    @("403×3=1209", "205×8=1640"),
    @("679×7=4753", "313×2=626"),
    @("348×2=696",  "380×3=1140"),
    @("639×6=3834", "287×2=574")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]

    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false,
                         $true, 1, $false, $new, 2) | Out-Null
}

$d.Save()
